# "Listo tif y rastro"
#
# The author finished working on the "Rastros(TIF)" sheet and moved over to
# the "Rastros (Rastro)" sheet: renamed it (dropped the stray space in the
# tab name), made it the active/selected tab, scrolled it rightwards so
# column Z is at the left edge of the view, selected cell AL4, and
# widened column AL (38) to fit its (now longer) header text.

$wb = $excel.ActiveWorkbook

# --- 1. Rename "Rastros (Rastro)" -> "Rastros(Rastro)" ----------------------
$wsTif    = $wb.Worksheets.Item(1)
$wsRastro = $wb.Worksheets.Item(2)
$wsRastro.Name = "Rastros(Rastro)"

# --- 2. Widen column AL (38) on the Rastro sheet to fit its header ----------
$wsRastro.Range("AL1").ColumnWidth = 17.5

# --- 3. Switch focus to the Rastro sheet: it becomes the active tab --------
$wsRastro.Activate()

# Scroll the view so column Z sits at the top-left and select AL4, mirroring
# where the user ended up after reviewing the widened column.
$excel.ActiveWindow.ScrollColumn = 26
$excel.ActiveWindow.ScrollRow = 1
$wsRastro.Range("AL4").Select() | Out-Null
